# Uncommented the code for cross browser/parallel testing.
# This resets the scratch "pythonCode" test-data sheet back to its
# minimal/blank fixture state (rows 4-11 removed, replaced by a single
# blank bordered row) and makes "pythonCode" the active sheet/tab again.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("pythonCode")

# Remove the old scratch python-code test rows (rows 4 through 11).
$ws.Rows("4:11").Delete()

# Add a blank, bordered row right after the remaining data (new row 12).
$newRow = $ws.Range("A12:C12")
$newRow.Borders.Weight = 2

# Update the selection shown on this sheet and make it the active tab.
$ws.Range("B6").Select()
$ws.Activate()
